$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.08104266666666667
$ws.Range("H2").Value = 0.243128
$ws.Range("I2").Value = 0.1109852216299026
$ws.Range("J2").Value = 0.1109852216299026
$ws.Range("O2").Value = 0.7533925587839766
$ws.Range("P2").Value = 0.7533925587839766
$ws.Range("Q2").Value = 0.001997701733333333
$ws.Range("R2").Value = 0.0179793156
$ws.Range("S2").Value = 0.08361544011095903
$ws.Range("T2").Value = 0.08361544011095905

# Row 3
$ws.Range("G3").Value = 0.08104266666666667
$ws.Range("H3").Value = 0.243128
$ws.Range("I3").Value = 0.1109852216299026
$ws.Range("J3").Value = 0.1109852216299026
$ws.Range("M3").Value = 0.008068666666666667
$ws.Range("N3").Value = 0.024206
$ws.Range("O3").Value = 0.2466074412160235
$ws.Range("P3").Value = 0.2466074412160235
$ws.Range("Q3").Value = 0.0006539062631111111
$ws.Range("R3").Value = 0.005885156368
$ws.Range("S3").Value = 0.02736978151894353
$ws.Range("T3").Value = 0.02736978151894353

# Row 4
$ws.Range("G4").Value = 0.5478883333333334
$ws.Range("I4").Value = 0.7503147490635131
$ws.Range("J4").Value = 0.7503147490635131
$ws.Range("O4").Value = 0.7533925587839766
$ws.Range("P4").Value = 0.7533925587839766
$ws.Range("S4").Value = 0.5652815486903174
$ws.Range("T4").Value = 0.5652815486903174

# Row 5
$ws.Range("G5").Value = 0.5478883333333334
$ws.Range("I5").Value = 0.7503147490635131
$ws.Range("J5").Value = 0.7503147490635131
$ws.Range("M5").Value = 0.008068666666666667
$ws.Range("N5").Value = 0.024206
$ws.Range("O5").Value = 0.2466074412160235
$ws.Range("P5").Value = 0.2466074412160235
$ws.Range("Q5").Value = 0.004420728332222223
$ws.Range("R5").Value = 0.03978655499
$ws.Range("S5").Value = 0.1850332003731957
$ws.Range("T5").Value = 0.1850332003731957

# Row 6
$ws.Range("G6").Value = 0.1012803333333333
$ws.Range("H6").Value = 0.303841
$ws.Range("I6").Value = 0.1387000293065843
$ws.Range("J6").Value = 0.1387000293065843
$ws.Range("O6").Value = 0.7533925587839766
$ws.Range("P6").Value = 0.7533925587839766
$ws.Range("Q6").Value = 0.002496560216666667
$ws.Range("R6").Value = 0.02246904195
$ws.Range("S6").Value = 0.1044955699827001
$ws.Range("T6").Value = 0.1044955699827001

# Row 7
$ws.Range("G7").Value = 0.1012803333333333
$ws.Range("H7").Value = 0.303841
$ws.Range("I7").Value = 0.1387000293065843
$ws.Range("J7").Value = 0.1387000293065843
$ws.Range("M7").Value = 0.008068666666666667
$ws.Range("N7").Value = 0.024206
$ws.Range("O7").Value = 0.2466074412160235
$ws.Range("P7").Value = 0.2466074412160235
$ws.Range("Q7").Value = 0.0008171972495555557
$ws.Range("R7").Value = 0.007354775246
$ws.Range("S7").Value = 0.03420445932388422
$ws.Range("T7").Value = 0.03420445932388422
